$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
#    We clone the existing bold/italic duplicate-title paragraph found near
#    the end of the doc (it already has the "<empty run> + <formatted run>"
#    shape we need) via Copy/Paste, then retarget its text - this preserves
#    the leading empty run that a from-scratch InsertAfter would collapse.
# ---------------------------------------------------------------------------
$title = $d.Paragraphs.Item(1)
$title.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# locate the duplicate bold title paragraph near the end of the document to
# use as a donor for the run/formatting shape (empty run + bold run); it is
# the second-to-last paragraph (the last one is the italic meta description)
$total = $d.Paragraphs.Count
$donor = $d.Paragraphs.Item($total - 1)
$donorLen = ($donor.Range.Text).Length
$donor.Range.Copy()
$metaPara.Range.Paste()

$metaPara = $d.Paragraphs.Item(2)
$metaStart = $metaPara.Range.Start
$boldLabel = "Meta description"
$rBold = $d.Range($metaStart, $metaStart + $donorLen)
$rBold.Text = $boldLabel

$metaPara = $d.Paragraphs.Item(2)
$metaEnd = $metaPara.Range.End
$restText = ": Read our review of Bugs Money and play for free. Exciting features and impeccable graphics, including Glow Wilds and Free Spins."
$rRest = $d.Range($metaEnd - 1, $metaEnd - 1)
$rRest.InsertAfter($restText)

# ---------------------------------------------------------------------------
# 2) Remove the duplicate bold "Play Bugs Money for Free..." paragraph that
#    used to sit right before the closing italic meta-description paragraph.
# ---------------------------------------------------------------------------
$total = $d.Paragraphs.Count
$dupPara = $d.Paragraphs.Item($total - 1)

$dupStart = $dupPara.Range.Start
$dupEnd = $dupPara.Range.End
$d.Range($dupStart, $dupEnd).Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the remaining (italic) closing paragraph with the
#    new image-prompt copy, keeping its italic run formatting intact.
# ---------------------------------------------------------------------------
$total = $d.Paragraphs.Count
$closingPara = $d.Paragraphs.Item($total)
$closingStart = $closingPara.Range.Start
$closingEnd = $closingPara.Range.End

$newClosingText = "Create a cartoon-style feature image for Bugs Money that prominently features a happy Maya warrior wearing glasses. The warrior should be standing in a bright green lawn with flowers and bugs surrounding them, holding a honeycomb with money flying out of it in the background. It should convey the fun and colorful nature of the game and showcase the potential for big wins. The image should be eye-catching and playful, enticing players to try their luck with Bugs Money."

$d.Range($closingStart, $closingEnd - 1).Text = $newClosingText

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
Write-Output "Paragraph 2: $($d.Paragraphs.Item(2).Range.Text)"
Write-Output "Last paragraph: $($d.Paragraphs.Item($d.Paragraphs.Count).Range.Text)"
